$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 6 de Septiembre de 2020 a las 14:36"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 6432209
$ws.Range("C4").Value = 1057
$ws.Range("D4").Value = 3707138
$ws.Range("E4").Value = 2532207
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 46
$ws.Range("H4").Value = 192864

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 4120846
$ws.Range("C6").Value = 10007
$ws.Range("D6").Value = 3183467
$ws.Range("E6").Value = 866647
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 53
$ws.Range("H6").Value = 70732

# Row 24: Alemania
$ws.Range("A24").Value = "Alemania"
$ws.Range("B24").Value = 251065
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 226208
$ws.Range("E24").Value = 15456
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 9401

# Row 40: Oman
$ws.Range("A40").Value = "Oman"
$ws.Range("B40").Value = 87072
$ws.Range("C40").Value = 692
$ws.Range("D40").Value = 82406
$ws.Range("E40").Value = 3938
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 23
$ws.Range("H40").Value = 728

# Row 44: Paises Bajos
$ws.Range("A44").Value = "Paises Bajos"
$ws.Range("B44").Value = 74787
$ws.Range("C44").Value = 925
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 6243

# Row 45: Emiratos Arabes Unidos
$ws.Range("A45").Value = "Emiratos Arabes Unidos"
$ws.Range("B45").Value = 73984
$ws.Range("C45").Value = 513
$ws.Range("D45").Value = 66095
$ws.Range("E45").Value = 7501
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 388

# Row 58: Nepal
$ws.Range("A58").Value = "Nepal"
$ws.Range("B58").Value = 46257
$ws.Range("C58").Value = 980
$ws.Range("D58").Value = 28941
$ws.Range("E58").Value = 17027
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 9
$ws.Range("H58").Value = 289

# Row 59: Argelia
$ws.Range("A59").Value = "Argelia"
$ws.Range("B59").Value = 46071
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 32481
$ws.Range("E59").Value = 12041
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 1549

# Row 75: Estado de Palestina
$ws.Range("A75").Value = "Estado de Palestina"
$ws.Range("B75").Value = 26127
$ws.Range("C75").Value = 552
$ws.Range("D75").Value = 16843
$ws.Range("E75").Value = 9103
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 181

# Row 82: Dinamarca
$ws.Range("A82").Value = "Dinamarca"
$ws.Range("B82").Value = 17883
$ws.Range("C82").Value = 147
$ws.Range("D82").Value = 15760
$ws.Range("E82").Value = 1496
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 627

# Row 83: Libia
$ws.Range("A83").Value = "Libia"
$ws.Range("B83").Value = 17749
$ws.Range("C83").Value = 655
$ws.Range("D83").Value = 2081
$ws.Range("E83").Value = 15383
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 13
$ws.Range("H83").Value = 285

# Row 128: Gambia
$ws.Range("A128").Value = "Gambia"
$ws.Range("B128").Value = 3197
$ws.Range("C128").Value = 46
$ws.Range("D128").Value = 1315
$ws.Range("E128").Value = 1783
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 99

# Row 129: Eslovenia
$ws.Range("A129").Value = "Eslovenia"
$ws.Range("B129").Value = 3165
$ws.Range("C129").Value = 43
$ws.Range("D129").Value = 2483
$ws.Range("E129").Value = 547
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 135

# Row 196: Curazao
$ws.Range("A196").Value = "Curazao"
$ws.Range("B196").Value = 92
$ws.Range("C196").Value = 4
$ws.Range("D196").Value = 45
$ws.Range("E196").Value = 46
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 1
